$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Point existing rows at the new "Staging_Env" subfolder for the expected file paths
$ws.Range("D2").Value = "\Testdata\Templates\ImportPublications\Staging_Env\Test dataset - Success Case Sheet.xlsx"
$ws.Range("D3").Value = "\Testdata\Templates\ImportPublications\Staging_Env\Test dataset - Failure Case Sheet.xlsx"

# Add new rows 4-7 with additional staging test data (populations pop3-pop6)
$ws.Range("A4").Value = "pop3"
$ws.Range("B4").Value = "Test - Test - Ovid search - 10/30/2020"
$ws.Range("C4").Value = "Test dataset - Header Mismatch.xlsx"
$ws.Range("D4").Value = "\Testdata\Templates\ImportPublications\Staging_Env\Test dataset - Header Mismatch.xlsx"
$ws.Range("A4").Style = "Normal"

$ws.Range("A5").Value = "pop4"
$ws.Range("B5").Value = "Test - Test - Ovid search - 10/30/2020"
$ws.Range("C5").Value = "Test dataset - Letters in Publication Identifier.xlsx"
$ws.Range("D5").Value = "\Testdata\Templates\ImportPublications\Staging_Env\Test dataset - Letters in Publication Identifier.xlsx"
$ws.Range("A5").Style = "Normal"

$ws.Range("A6").Value = "pop5"
$ws.Range("B6").Value = "Test - Test - Ovid search - 10/30/2020"
$ws.Range("C6").Value = "Test dataset - Empty value in Publication Identifier.xlsx"
$ws.Range("D6").Value = "\Testdata\Templates\ImportPublications\Staging_Env\Test dataset - Empty value in Publication Identifier.xlsx"
$ws.Range("A6").Style = "Normal"

$ws.Range("A7").Value = "pop6"
$ws.Range("B7").Value = "Test - Test - Ovid search - 10/30/2020"
$ws.Range("C7").Value = "Test dataset - Duplicate value in FA18 column.xlsx"
$ws.Range("D7").Value = "\Testdata\Templates\ImportPublications\Staging_Env\Test dataset - Duplicate value in FA18 column.xlsx"
$ws.Range("A7").Style = "Normal"

# Column C widened slightly to fit the new, longer file names (bestFit target ~32.21875)
$ws.Columns("C").ColumnWidth = 31.39

# Final selection in the saved workbook
$ws.Range("C7").Select()
